$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text corrections / reshuffle in column C (disease names) ---
$ws.Range("C6").Value  = "Common rust_"
$ws.Range("C13").Value = "Early Blight"
$ws.Range("C14").Value = "Late Blight"
$ws.Range("C15").Value = "Leaf Mold"
$ws.Range("C16").Value = "Septoria Leaf Spot"
$ws.Range("C17").Value = "Two-spotted Spider Mites"
$ws.Range("C18").Value = "Target Spot"
$ws.Range("C19").Value = "Mosaic virus"
$ws.Range("C20").Value = "Yellow Leaf Curl Virus"
$ws.Range("C21").Value = "Healthy"
$ws.Range("C22").Value = "Bacterial spot"

# --- Fill in the "Jumlah" (count) column D with the sample counts ---
$ws.Range("D5").Value  = 410
$ws.Range("D6").Value  = 500
$ws.Range("D7").Value  = 500
$ws.Range("D8").Value  = 500
$ws.Range("D9").Value  = 500
$ws.Range("D10").Value = 500
$ws.Range("D11").Value = 121
$ws.Range("D12").Value = 500
$ws.Range("D13").Value = 500
$ws.Range("D14").Value = 500
$ws.Range("D15").Value = 500
$ws.Range("D16").Value = 500
$ws.Range("D17").Value = 500
$ws.Range("D18").Value = 500
$ws.Range("D19").Value = 299
$ws.Range("D20").Value = 500
$ws.Range("D21").Value = 500
$ws.Range("D22").Value = 500

# Center the new counts both horizontally and vertically, matching the
# rest of the table's look.
$ws.Range("D5:D22").HorizontalAlignment = -4108
$ws.Range("D5:D22").VerticalAlignment = -4108

# "Testing" header gets centered too.
$ws.Range("E4").HorizontalAlignment = -4108

# --- Grand total row ---
$ws.Range("D23").Formula = "=SUM(D5:D22)"

$ws.Range("D23").Select()

$wb.Application.CalculateFull()
